$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -----------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- Formatting (bold font, thin box border, centered/top aligned) ----
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.Borders.LineStyle = 1
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4160

# Copy the finished format onto A2 so both cells resolve to the very
# same style record instead of each incrementally growing its own.
$b1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
